$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1: extend the function-index header from 10 columns (A:J) out to
# 80 columns (A:CB). Columns A:J already hold 1,1,2,2,3,3,4,4,5,5 (functions
# 1-5, one column per bound). Continue the same "two columns per function"
# pattern for functions 6-40 across columns K:CB.
$col = 11
$funcIndex = 6
while ($col -le 80) {
    $ws.Cells.Item(1, $col).Value = $funcIndex
    $ws.Cells.Item(1, $col + 1).Value = $funcIndex
    $col = $col + 2
    $funcIndex = $funcIndex + 1
}

# --- Rows 2-11: the old bound values lived in columns C:J (functions 2-5).
# Clear them out - the new layout places bounds for functions 19, 28 and 33
# (and a couple of extra rows for function 37) further to the right so the
# sheet can hold bounds for all 40 functions, with only f1-f5 runnable.
$ws.Range("C2:J11").ClearContents()

for ($row = 2; $row -le 11; $row++) {
    # function 1 (already in place) - lower/upper bound
    $ws.Cells.Item($row, 1).Value = -32
    $ws.Cells.Item($row, 2).Value = 32

    # function 19 -> columns AK:AL
    $ws.Cells.Item($row, 37).Value = -30
    $ws.Cells.Item($row, 38).Value = 30

    # function 28 -> columns BC:BD
    $ws.Cells.Item($row, 55).Value = -600
    $ws.Cells.Item($row, 56).Value = 600

    # function 33 -> columns BM:BN
    $ws.Cells.Item($row, 65).Value = -10
    $ws.Cells.Item($row, 66).Value = 10
}

# function 37 -> columns BU:BV, only populated on rows 2 and 3
$ws.Cells.Item(2, 73).Value = 0
$ws.Cells.Item(2, 74).Value = 1
$ws.Cells.Item(3, 73).Value = 0
$ws.Cells.Item(3, 74).Value = 1

# --- Update the view: scroll to show the newly added columns and move the
# active selection to where the author last left it.
$ws.Application.GoTo($ws.Range("BS10"), $true)
$ws.Range("BS10").Select()
